# Update the "取得日時" (acquisition timestamp) column on the ランサーズ sheet
# for the existing data rows (2-12) to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-07 18:29:25"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
